$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1196.809
$ws.Range("I15").Value = 1196.809
$ws.Range("K15").Value = 3590.427
$ws.Range("M15").Value = -3421.427

$ws.Range("H41").Value = 33337424
$ws.Range("I41").Value = 555.36365
$ws.Range("K41").Value = 555.36365
$ws.Range("M41").Value = -115.36365

$ws.Range("H88").Value = 1625.25
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 1625.25
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 1625.25
$ws.Range("N88").Value = -2437.25
$ws.Range("M88").ClearContents()

$ws.Range("H91").Value = 1625.25
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 1625.25
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 1625.25
$ws.Range("N91").Value = -4433.25
$ws.Range("M91").ClearContents()

$ws.Range("H125").Value = 12350323
$ws.Range("I125").Value = 794
$ws.Range("J125").Value = 13894014
$ws.Range("K125").Value = 7146
$ws.Range("L125").Value = 125046126
$ws.Range("M125").Value = -4686
$ws.Range("N125").Value = -125051046

$ws.Range("H129").Value = 2425.7273
$ws.Range("I129").Value = 1048.5
$ws.Range("J129").Value = 3212.7144
$ws.Range("K129").Value = 3145.5
$ws.Range("L129").Value = 9638.143199999999
$ws.Range("M129").Value = 1854.5
$ws.Range("N129").Value = -19638.1432

$ws.Range("H132").Value = 1856.8334
$ws.Range("I132").Value = 1780.25
$ws.Range("J132").Value = 2699.25
$ws.Range("K132").Value = 5340.75
$ws.Range("L132").Value = 8097.75
$ws.Range("M132").Value = -2810.75
$ws.Range("N132").Value = -13157.75

$ws.Range("H135").Value = 1932.2222
$ws.Range("I135").Value = 1547.2
$ws.Range("K135").Value = 13924.8
$ws.Range("M135").Value = -11389.8

$ws.Range("H138").Value = 6643.5947
$ws.Range("J138").Value = 7578.263
$ws.Range("L138").Value = 22734.789
$ws.Range("N138").Value = -33014.789

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5337.5674
$ws.Range("I32").Value = 3984.5454
$ws.Range("J32").Value = 16500
$ws.Range("K32").Value = 3984.5454
$ws.Range("L32").Value = 16500
$ws.Range("M32").Value = -3697.5454
$ws.Range("N32").Value = -17074

$ws.Range("H74").Value = 2072.7104
$ws.Range("I74").Value = 2140.0881
$ws.Range("J74").Value = 1500
$ws.Range("K74").Value = 2140.0881
$ws.Range("L74").Value = 1500
$ws.Range("M74").Value = -1266.0881
$ws.Range("N74").Value = -3248

$ws.Range("H77").Value = 2072.7104
$ws.Range("I77").Value = 2140.0881
$ws.Range("J77").Value = 1500
$ws.Range("K77").Value = 10700.4405
$ws.Range("L77").Value = 7500
$ws.Range("M77").Value = -6332.440500000001
$ws.Range("N77").Value = -16236

$ws.Range("H102").Value = 1705.6471
$ws.Range("I102").Value = 1682.7333
$ws.Range("K102").Value = 1682.7333
$ws.Range("M102").Value = -60.7333000000001

$ws.Range("H124").Value = 46875
$ws.Range("J124").Value = 46875
$ws.Range("L124").Value = 46875
$ws.Range("N124").Value = -56695

$ws.Range("H125").Value = 92507.71000000001
$ws.Range("J125").Value = 92507.71000000001
$ws.Range("L125").Value = 92507.71000000001
$ws.Range("N125").Value = -102347.71

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 711908.0600000001
$ws.Range("I86").Value = 898252.8
$ws.Range("J86").Value = 3798
$ws.Range("K86").Value = 898252.8
$ws.Range("L86").Value = 3798
$ws.Range("M86").Value = -897129.8
$ws.Range("N86").Value = -6044

$ws.Range("H89").Value = 711908.0600000001
$ws.Range("I89").Value = 898252.8
$ws.Range("J89").Value = 3798
$ws.Range("K89").Value = 4491264
$ws.Range("L89").Value = 18990
$ws.Range("M89").Value = -4485648
$ws.Range("N89").Value = -30222

$ws.Range("H134").Value = 23067.08
$ws.Range("I134").Value = 3169.2341
$ws.Range("J134").Value = 334800
$ws.Range("K134").Value = 9507.702300000001
$ws.Range("L134").Value = 1004400
$ws.Range("M134").Value = -6972.702300000001
$ws.Range("N134").Value = -1009470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 5995.5
$ws.Range("I25").Value = 5995.5
$ws.Range("K25").Value = 5995.5
$ws.Range("M25").Value = -5821.5

$ws.Range("H31").Value = 50220.383
$ws.Range("I31").Value = 1839
$ws.Range("J31").Value = 103439.9
$ws.Range("K31").Value = 1839
$ws.Range("L31").Value = 103439.9
$ws.Range("M31").Value = -1544
$ws.Range("N31").Value = -104029.9

$ws.Range("H34").Value = 50220.383
$ws.Range("I34").Value = 1839
$ws.Range("J34").Value = 103439.9
$ws.Range("K34").Value = 1839
$ws.Range("L34").Value = 103439.9
$ws.Range("M34").Value = -1637
$ws.Range("N34").Value = -103843.9

$ws.Range("H58").Value = 2112.04
$ws.Range("I58").Value = 1945.5
$ws.Range("K58").Value = 1945.5
$ws.Range("M58").Value = -1742.5

$ws.Range("H94").Value = 1448.3334
$ws.Range("I94").Value = 1965
$ws.Range("J94").Value = 931.6667
$ws.Range("K94").Value = 1965
$ws.Range("L94").Value = 931.6667
$ws.Range("M94").Value = -1514
$ws.Range("N94").Value = -1833.6667

$ws.Range("H107").Value = 340.75
$ws.Range("I107").Value = 285.66666
$ws.Range("J107").Value = 506
$ws.Range("K107").Value = 285.66666
$ws.Range("L107").Value = 506
$ws.Range("M107").Value = 1634.33334
$ws.Range("N107").Value = -4346

$ws.Range("H123").Value = 63945
$ws.Range("J123").Value = 63945
$ws.Range("L123").Value = 63945
$ws.Range("N123").Value = -73745

$ws.Range("H130").Value = 77267
$ws.Range("J130").Value = 77267
$ws.Range("L130").Value = 77267
$ws.Range("N130").Value = -87307

$ws.Range("H136").Value = 2112.04
$ws.Range("I136").Value = 1945.5
$ws.Range("K136").Value = 5836.5
$ws.Range("M136").Value = -3286.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 3140.625
$ws.Range("J98").Value = 3002.5557
$ws.Range("L98").Value = 9007.667099999999
$ws.Range("N98").Value = -12003.6671

$ws.Range("H106").Value = 38427.215
$ws.Range("J106").Value = 39152.383
$ws.Range("L106").Value = 117457.149
$ws.Range("N106").Value = -119349.149

$ws.Range("H107").Value = 90235.05
$ws.Range("I107").Value = 1144.125
$ws.Range("J107").Value = 137750.2
$ws.Range("K107").Value = 3432.375
$ws.Range("L107").Value = 413250.6
$ws.Range("M107").Value = -1512.375
$ws.Range("N107").Value = -417090.6

$ws.Range("H139").Value = 5699.558
$ws.Range("I139").Value = 2932.889
$ws.Range("K139").Value = 8798.667000000001
$ws.Range("M139").Value = -3658.667000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 151108.78
$ws.Range("I132").Value = 9593.091
$ws.Range("K132").Value = 28779.273
$ws.Range("M132").Value = -26249.273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 115530.664
$ws.Range("I40").Value = 201596.6
$ws.Range("K40").Value = 201596.6
$ws.Range("M40").Value = -201460.6

$ws.Range("H136").Value = 678436.2
$ws.Range("I136").Value = 2007809.6
$ws.Range("J136").Value = 13749.5
$ws.Range("K136").Value = 6023428.800000001
$ws.Range("L136").Value = 41248.5
$ws.Range("M136").Value = -6020878.800000001
$ws.Range("N136").Value = -46348.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 38345
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H132").Value = 167315.67
$ws.Range("I132").Value = 723.5
$ws.Range("K132").Value = 2170.5
$ws.Range("M132").Value = 359.5

$ws.Range("H136").Value = 356590.25
$ws.Range("I136").Value = 374078.28
$ws.Range("K136").Value = 1122234.84
$ws.Range("M136").Value = -1119684.84
